# Re-create the "Updated cryptos list" GitHub Actions commit: refresh the
# Price (D) and Volume(1h) (E) columns for every coin row, and swap the
# Mantle/Stellar rows (48/49) which changed order upstream.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.383.81"
$ws.Range("E2").Value = "  -2.53%  "

$ws.Range("D3").Value = "2.379.64"
$ws.Range("E3").Value = "  -2.96%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.77%  "

$ws.Range("E8").Value = "  -1.49%  "

$ws.Range("D9").Value = "2.401.79"
$ws.Range("E9").Value = "  -1.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0961"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.78%  "

$ws.Range("E11").Value = "  -0.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.317"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.96%  "

$ws.Range("D14").Value = "2.799.90"
$ws.Range("E14").Value = "  -3.22%  "

$ws.Range("D15").Value = "56.251.02"
$ws.Range("E15").Value = "  -2.65%  "

$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("D18").Value = "2.397.25"
$ws.Range("E18").Value = "  -2.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "310.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("D27").Value = "2.488.73"
$ws.Range("E27").Value = "  -3.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.376"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.37%  "

$ws.Range("E29").Value = "  -4.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.37%  "

$ws.Range("D32").Value = "0.0₃0718"
$ws.Range("E32").Value = "  -2.09%  "

$ws.Range("E33").Value = "  -1.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.29%  "

$ws.Range("E35").Value = "  -3.63%  "

$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.60%  "

$ws.Range("E39").Value = "  +2.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.801"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "129.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "252.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.44%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0907"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.49%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.561"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0487"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.60%  "
